$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 8).Value = $false
}

$ws.Range("K13").Select()
